$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.798.34"
$ws.Range("E2").Value = "  +0.20%  "

# Row 3
$ws.Range("D3").Value = "3.810.12"
$ws.Range("E3").Value = "  +0.61%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").Value = "'602.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.22%  "

# Row 6
$ws.Range("D6").Value = "'166.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.55%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").Value = "'0.518"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.09%  "

# Row 9
$ws.Range("E9").Value = "  +0.21%  "

# Row 10
$ws.Range("E10").Value = "  +0.85%  "

# Row 11
$ws.Range("D11").Value = "'6.45"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.09%  "

# Row 12
$ws.Range("D12").Value = "'0.0000250"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.06%  "

# Row 13
$ws.Range("D13").Value = "'36.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.23%  "

# Row 14
$ws.Range("D14").Value = "4.447.03"
$ws.Range("E14").Value = "  +0.57%  "

# Row 15
$ws.Range("D15").Value = "3.798.38"
$ws.Range("E15").Value = "  +0.33%  "

# Row 16
$ws.Range("D16").Value = "67.809.35"
$ws.Range("E16").Value = "  +0.26%  "

# Row 17
$ws.Range("D17").Value = "'18.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.50%  "

# Row 18
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'7.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.78%  "

# Row 19
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "'0.113"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.93%  "

# Row 20
$ws.Range("D20").Value = "'463.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.04%  "

# Row 21
$ws.Range("D21").Value = "'9.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.93%  "

# Row 22
$ws.Range("E22").Value = "  +1.03%  "

# Row 23
$ws.Range("D23").Value = "'0.0000147"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.64%  "

# Row 24
$ws.Range("D24").Value = "'83.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.17%  "

# Row 25
$ws.Range("D25").Value = "'12.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.73%  "

# Row 26
$ws.Range("E26").Value = "  -0.09%  "

# Row 27
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.05%  "

# Row 28
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'10.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.01%  "

# Row 29
$ws.Range("D29").Value = "3.958.56"
$ws.Range("E29").Value = "  +0.65%  "

# Row 30
$ws.Range("E30").Value = "  +0.11%  "

# Row 31
$ws.Range("D31").Value = "'7.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.23%  "

# Row 32
$ws.Range("D32").Value = "'2.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.03%  "

# Row 33
$ws.Range("D33").Value = "'29.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.66%  "

# Row 34
$ws.Range("E34").Value = "  +0.14%  "

# Row 35
$ws.Range("E35").Value = "  -0.18%  "

# Row 36
$ws.Range("D36").Value = "'0.100"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.13%  "

# Row 37
$ws.Range("E37").Value = "  +0.07%  "

# Row 38
$ws.Range("D38").Value = "'0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.45%  "

# Row 39
$ws.Range("D39").Value = "'5.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.95%  "

# Row 40
$ws.Range("D40").Value = "'3.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.50%  "

# Row 41
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.02%  "

# Row 42
$ws.Range("E42").Value = "  -0.01%  "

# Row 43
$ws.Range("D43").Value = "'44.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.02%  "

# Row 44
$ws.Range("D44").Value = "'47.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.92%  "

# Row 45
$ws.Range("E45").Value = "  +0.07%  "

# Row 46
$ws.Range("D46").Value = "'27.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.48%  "

# Row 47
$ws.Range("D47").Value = "'151.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.13%  "

# Row 48
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "'8.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.55%  "

# Row 49
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "'1.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.96%  "

# Row 50
$ws.Range("D50").Value = "'1.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.23%  "

# Row 51
$ws.Range("D51").Value = "'390.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.86%  "
